$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 432.4
$ws.Range("I55").Value = 662
$ws.Range("J55").Value = 202.8
$ws.Range("K55").Value = 662
$ws.Range("L55").Value = 202.8
$ws.Range("M55").Value = -448
$ws.Range("N55").Value = -630.8
$ws.Range("H64").Value = 4526.087
$ws.Range("I64").Value = 4975
$ws.Range("J64").Value = 3500
$ws.Range("K64").Value = 4975
$ws.Range("L64").Value = 3500
$ws.Range("M64").Value = -4727
$ws.Range("N64").Value = -3996
$ws.Range("H67").Value = 4526.087
$ws.Range("I67").Value = 4975
$ws.Range("J67").Value = 3500
$ws.Range("K67").Value = 4975
$ws.Range("L67").Value = 3500
$ws.Range("M67").Value = -4117
$ws.Range("N67").Value = -5216
$ws.Range("H92").Value = 86806110
$ws.Range("I92").Value = 3968891
$ws.Range("K92").Value = 3968891
$ws.Range("M92").Value = -3967643
$ws.Range("H129").Value = 1068.809
$ws.Range("I129").Value = 830
$ws.Range("J129").Value = 1086.0723
$ws.Range("K129").Value = 2490
$ws.Range("L129").Value = 3258.2169
$ws.Range("M129").Value = 2510
$ws.Range("N129").Value = -13258.2169
$ws.Range("H132").Value = 2212.182
$ws.Range("I132").Value = 1845.7894
$ws.Range("J132").Value = 4532.6665
$ws.Range("K132").Value = 5537.3682
$ws.Range("L132").Value = 13597.9995
$ws.Range("M132").Value = -3007.3682
$ws.Range("N132").Value = -18657.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 673.27905
$ws.Range("I2").Value = 535.8108
$ws.Range("J2").Value = 1521
$ws.Range("K2").Value = 535.8108
$ws.Range("L2").Value = 1521
$ws.Range("M2").Value = -422.8108
$ws.Range("N2").Value = -1747
$ws.Range("H32").Value = 2330.55
$ws.Range("I32").Value = 2198.4639
$ws.Range("K32").Value = 2198.4639
$ws.Range("M32").Value = -1911.4639
$ws.Range("H74").Value = 1465.6875
$ws.Range("I74").Value = 1063.9062
$ws.Range("K74").Value = 1063.9062
$ws.Range("M74").Value = -189.9061999999999
$ws.Range("H77").Value = 1465.6875
$ws.Range("I77").Value = 1063.9062
$ws.Range("K77").Value = 5319.530999999999
$ws.Range("M77").Value = -951.530999999999
$ws.Range("H116").Value = 673.27905
$ws.Range("I116").Value = 535.8108
$ws.Range("J116").Value = 1521
$ws.Range("K116").Value = 535.8108
$ws.Range("L116").Value = 1521
$ws.Range("M116").Value = 1758.1892
$ws.Range("N116").Value = -6109
$ws.Range("H122").Value = 952216.9
$ws.Range("I122").Value = 1223690.2
$ws.Range("K122").Value = 3671070.6
$ws.Range("M122").Value = -3668620.6
$ws.Range("H132").Value = 2635573.5
$ws.Range("I132").Value = 3006.913
$ws.Range("J132").Value = 6672175.5
$ws.Range("K132").Value = 9020.739
$ws.Range("L132").Value = 20016526.5
$ws.Range("M132").Value = -6490.739
$ws.Range("N132").Value = -20021586.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 673.27905
$ws.Range("I3").Value = 535.8108
$ws.Range("J3").Value = 1521
$ws.Range("K3").Value = 535.8108
$ws.Range("L3").Value = 1521
$ws.Range("M3").Value = -421.8108
$ws.Range("N3").Value = -1749
$ws.Range("H122").Value = 70000
$ws.Range("J122").Value = 70000
$ws.Range("L122").Value = 70000
$ws.Range("N122").Value = -79800
$ws.Range("H135").Value = 78926.664
$ws.Range("J135").Value = 78926.664
$ws.Range("L135").Value = 78926.664
$ws.Range("N135").Value = -89066.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2116.647
$ws.Range("I16").Value = 1776.6666
$ws.Range("J16").Value = 2499.125
$ws.Range("K16").Value = 1776.6666
$ws.Range("L16").Value = 2499.125
$ws.Range("M16").Value = -1489.6666
$ws.Range("N16").Value = -3073.125
$ws.Range("H22").Value = 1039.3572
$ws.Range("I22").Value = 946.375
$ws.Range("J22").Value = 1163.3334
$ws.Range("K22").Value = 946.375
$ws.Range("L22").Value = 1163.3334
$ws.Range("M22").Value = -596.375
$ws.Range("N22").Value = -1863.3334
$ws.Range("H31").Value = 4823.3335
$ws.Range("I31").Value = 1944.6052
$ws.Range("J31").Value = 10580.789
$ws.Range("K31").Value = 1944.6052
$ws.Range("L31").Value = 10580.789
$ws.Range("M31").Value = -1649.6052
$ws.Range("N31").Value = -11170.789
$ws.Range("H34").Value = 4823.3335
$ws.Range("I34").Value = 1944.6052
$ws.Range("J34").Value = 10580.789
$ws.Range("K34").Value = 1944.6052
$ws.Range("L34").Value = 10580.789
$ws.Range("M34").Value = -1742.6052
$ws.Range("N34").Value = -10984.789
$ws.Range("H113").Value = 2116.647
$ws.Range("I113").Value = 1776.6666
$ws.Range("J113").Value = 2499.125
$ws.Range("K113").Value = 1776.6666
$ws.Range("L113").Value = 2499.125
$ws.Range("M113").Value = 393.3334
$ws.Range("N113").Value = -6839.125
$ws.Range("H132").Value = 2020.9412
$ws.Range("I132").Value = 1638.1765
$ws.Range("J132").Value = 2403.7058
$ws.Range("K132").Value = 4914.529500000001
$ws.Range("L132").Value = 7211.117400000001
$ws.Range("M132").Value = -2384.529500000001
$ws.Range("N132").Value = -12271.1174
$ws.Range("H134").Value = 315785.97
$ws.Range("I134").Value = 3372.0386
$ws.Range("K134").Value = 10116.1158
$ws.Range("M134").Value = -7581.1158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 617.64703
$ws.Range("I2").Value = 693.86664
$ws.Range("J2").Value = 46
$ws.Range("K2").Value = 4163.199839999999
$ws.Range("L2").Value = 276
$ws.Range("M2").Value = -4050.199839999999
$ws.Range("N2").Value = -502
$ws.Range("H38").Value = 4000326
$ws.Range("I38").Value = 5882661.5
$ws.Range("J38").Value = 362.75
$ws.Range("K38").Value = 17647984.5
$ws.Range("L38").Value = 1088.25
$ws.Range("M38").Value = -17647637.5
$ws.Range("N38").Value = -1782.25
$ws.Range("H107").Value = 347.86206
$ws.Range("I107").Value = 244.41667
$ws.Range("J107").Value = 420.88235
$ws.Range("K107").Value = 733.25001
$ws.Range("L107").Value = 1262.64705
$ws.Range("M107").Value = 1186.74999
$ws.Range("N107").Value = -5102.64705
$ws.Range("H131").Value = 1695841.1
$ws.Range("I131").Value = 5882878
$ws.Range("J131").Value = 1087.9762
$ws.Range("K131").Value = 17648634
$ws.Range("L131").Value = 3263.9286
$ws.Range("M131").Value = -17643594
$ws.Range("N131").Value = -13343.9286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 57167764
$ws.Range("I122").Value = 62637870
$ws.Range("K122").Value = 187913610
$ws.Range("M122").Value = -187911160

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5406.375
$ws.Range("J2").Value = 5406.375
$ws.Range("L2").Value = 5406.375
$ws.Range("N2").Value = -5630.375
$ws.Range("H132").Value = 9809433
$ws.Range("I132").Value = 10758281
$ws.Range("K132").Value = 32274843
$ws.Range("M132").Value = -32272313

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2106.577
$ws.Range("I132").Value = 1419.6316
$ws.Range("J132").Value = 3971.1428
$ws.Range("K132").Value = 4258.8948
$ws.Range("L132").Value = 11913.4284
$ws.Range("M132").Value = -1728.8948
$ws.Range("N132").Value = -16973.4284
$ws.Range("H136").Value = 2372.6736
$ws.Range("I136").Value = 2357.9092
$ws.Range("K136").Value = 7073.7276
$ws.Range("M136").Value = -4523.7276
